# Weekly price-sheet update: a new week (2023-04-05 / serial 45021) of
# "Frutilla" price data is inserted as a new row right before the existing
# row 379, shifting all subsequent rows down by one (old row 379 -> 380,
# ... old row 420 -> 421). The new row's Market/Region/Category columns
# (A,B,C,E,F,G,H,I,J,K,T) carry the same constant values used throughout
# this sheet; L-S mirror the most recent prior week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 379; everything from 379 down shifts to 380+.
$ws.Rows.Item(379).Insert()

$ws.Range("A379").Value = 4
$ws.Range("B379").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C379").Value = "Los Lagos"
$ws.Range("D379").Value = 45021
$ws.Range("E379").Value = 10
$ws.Range("F379").Value = "Fruta"
$ws.Range("G379").Value = 100101
$ws.Range("H379").Value = "Berries"
$ws.Range("I379").Value = 100112025
$ws.Range("J379").Value = "Frutilla"
$ws.Range("K379").Value = "Sin especificar"
$ws.Range("L379").Value = "Primera"
$ws.Range("M379").Value = 200
$ws.Range("N379").Value = 9500
$ws.Range("O379").Value = 10000
$ws.Range("P379").Value = 9750
$ws.Range("Q379").Value = "$/caja 7 kilos"
$ws.Range("R379").Value = "Región de La Araucanía"
$ws.Range("S379").Value = 1393
$ws.Range("T379").Value = 7
